$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full 3x3 cluster combination (ECs, FAPs, sCs) x (ECs, FAPs, sCs) for the
# Sfrp1 -> Fzd6 ligand-receptor pair, replacing the previous partial table.
$data = @(
    @("ECs","Sfrp1","Fzd6","ECs",2,0.6666666666666666,0.217473,0.652419,0.003819775075312922,0.003819775075312921,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,2.700135199188,24.301216792692,0.003284258772512135,0.003284258772512135),
    @("ECs","Sfrp1","Fzd6","FAPs",2,0.6666666666666666,0.217473,0.652419,0.003819775075312922,0.003819775075312921,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,0.3796816887489999,3.417135198741,0.0004618186961197803,0.0004618186961197802),
    @("ECs","Sfrp1","Fzd6","sCs",2,0.6666666666666666,0.217473,0.652419,0.003819775075312922,0.003819775075312921,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,0.060590080039,0.5453107203509999,0.00007369760668100643,0.00007369760668100641),
    @("FAPs","Sfrp1","Fzd6","ECs",3,1,55.43187966666667,166.295639,0.9736257481548445,0.9736257481548444,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,688.239778938628,6194.158010447652,0.837127538002819,0.8371275380028189),
    @("FAPs","Sfrp1","Fzd6","FAPs",3,1,55.43187966666667,166.295639,0.9736257481548445,0.9736257481548444,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,96.77739159514677,870.9965243563209,0.1177133639170313,0.1177133639170313),
    @("FAPs","Sfrp1","Fzd6","sCs",3,1,55.43187966666667,166.295639,0.9736257481548445,0.9736257481548444,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,15.44385751663678,138.994717649731,0.01878484623499412,0.01878484623499413),
    @("sCs","Sfrp1","Fzd6","ECs",3,1,1.284104333333334,3.852313000000001,0.02255447676984262,0.02255447676984261,3,1,12.415956,37.247868,0.8598042313376485,0.8598042313376485,15.943382902076,143.490446118684,0.01939243456231738,0.01939243456231738),
    @("sCs","Sfrp1","Fzd6","FAPs",3,1,1.284104333333334,3.852313000000001,0.02255447676984262,0.02255447676984261,3,1,1.745879666666666,5.237639,0.1209020654395331,0.1209020654395331,2.241891645445222,20.177024809007,0.002726882826381941,0.00272688282638194),
    @("sCs","Sfrp1","Fzd6","sCs",3,1,1.284104333333334,3.852313000000001,0.02255447676984262,0.02255447676984261,2,0.6666666666666666,0.2786096666666666,0.8358289999999999,0.01929370322281843,0.01929370322281843,0.3577638802752223,3.219874922477,0.0004351593811432958,0.0004351593811432958)
)

$row = 2
foreach ($rec in $data) {
    for ($col = 1; $col -le $rec.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rec[$col - 1]
    }
    $row = $row + 1
}
